# Fix Training Data Issue (#48)
# The Date column (BF) values were off by one day due to the way NBA
# stats were shown (e.g. "6-18-2012-13"). Correct them to the actual
# ISO-formatted game date "2013-06-18" for every data row (BF2:BF31),
# while keeping the values stored as plain text (not auto-converted to
# an Excel date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")

# Force text interpretation so Excel doesn't reinterpret the
# "YYYY-MM-DD"-looking string as a date serial number, then drop the
# temporary formatting again so the cells keep their original (default)
# style.
$rng.NumberFormat = "@"
$rng.Value = "2013-06-18"
$rng.ClearFormats()
